$d = $word.ActiveDocument
$cr = [char]13

# Finds a paragraph whose visible text (paragraph mark excluded) exactly
# matches $text. Searches the whole document story.
function Get-ParaByExactText {
    param($doc, $text)
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -eq ($text + $cr)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) New vocabulary line "Stornieren = Annuler" right after
#    "Bestellen = Commander".
# ---------------------------------------------------------------------
$i = 1
$bestellenIdx = -1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq ("Bestellen = Commander" + $cr)) {
        $bestellenIdx = $i
    }
    $i = $i + 1
}

if ($bestellenIdx -gt 0) {
    $pBestellen = $d.Paragraphs($bestellenIdx)
    [void]$pBestellen.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs($bestellenIdx + 1)
    $newPara.Range.Text = "Stornieren = Annuler"
}

# ---------------------------------------------------------------------
# 2) Re-rendered pagination markers (<w:lastRenderedPageBreak/>) shift
#    one line earlier in two spots. Surgically rewrite the affected
#    paragraphs (preserving their ids / run formatting) using InsertXML.
# ---------------------------------------------------------------------
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# "Signature:" gains the page-break marker.
$pSignature = Get-ParaByExactText $d "Signature:"
if ($pSignature -ne $null) {
    $xml = '<w:p ' + $wNs + ' w14:paraId="35C4174B" w14:textId="77777777" w:rsidR="00892EA5" w:rsidRPr="00892EA5" w:rsidRDefault="00892EA5"><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr></w:pPr><w:r w:rsidRPr="00892EA5"><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:lastRenderedPageBreak/><w:t>Signature:</w:t></w:r></w:p>'
    [void]$pSignature.Range.InsertXML($xml)
}

# "Si Minuscule  tutoyer" loses the page-break marker.
$pSiMin = Get-ParaByExactText $d "Si Minuscule  tutoyer"
if ($pSiMin -ne $null) {
    $xml = '<w:p ' + $wNs + ' w14:paraId="5AE82232" w14:textId="7AFAE37B" w:rsidR="00892EA5" w:rsidRDefault="00892EA5"><w:pPr><w:rPr><w:lang w:val="de-AT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-AT"/></w:rPr><w:t xml:space="preserve">Si Minuscule </w:t></w:r><w:r w:rsidRPr="00892EA5"><w:rPr><w:lang w:val="de-AT"/></w:rPr><w:sym w:font="Wingdings" w:char="F0E0"/></w:r><w:r><w:rPr><w:lang w:val="de-AT"/></w:rPr><w:t xml:space="preserve"> tutoyer</w:t></w:r></w:p>'
    [void]$pSiMin.Range.InsertXML($xml)
}

# "Autre signature:" loses the page-break marker (no replacement spot).
$pAutre = Get-ParaByExactText $d "Autre signature:"
if ($pAutre -ne $null) {
    $xml = '<w:p ' + $wNs + ' w14:paraId="2BF38BC2" w14:textId="2D46B5D7" w:rsidR="00892EA5" w:rsidRPr="00892EA5" w:rsidRDefault="00892EA5"><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr><w:r><w:rPr><w:u w:val="single"/></w:rPr><w:t>A</w:t></w:r><w:r w:rsidRPr="00892EA5"><w:rPr><w:u w:val="single"/></w:rPr><w:t>utre signature:</w:t></w:r></w:p>'
    [void]$pAutre.Range.InsertXML($xml)
}

# "Die Sitten = les coutumes" gains the page-break marker.
$pSitten = Get-ParaByExactText $d "Die Sitten = les coutumes"
if ($pSitten -ne $null) {
    $xml = '<w:p ' + $wNs + ' w14:paraId="143F6D72" w14:textId="618E1CCC" w:rsidR="00B67151" w:rsidRDefault="00B67151"><w:pPr><w:rPr><w:lang w:val="de-AT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-AT"/></w:rPr><w:lastRenderedPageBreak/><w:t>Die Sitten = les coutumes</w:t></w:r></w:p>'
    [void]$pSitten.Range.InsertXML($xml)
}

# "Die Bildung = la formation" loses the page-break marker.
$pBildung = Get-ParaByExactText $d "Die Bildung = la formation"
if ($pBildung -ne $null) {
    $xml = '<w:p ' + $wNs + ' w14:paraId="1AD1144E" w14:textId="5789A33C" w:rsidR="001E3499" w:rsidRDefault="00B67151" w:rsidP="00B67151"><w:pPr><w:rPr><w:lang w:val="de-AT"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="de-AT"/></w:rPr><w:t>Die Bildung = la formation</w:t></w:r></w:p>'
    [void]$pBildung.Range.InsertXML($xml)
}

Write-Output "edit applied"
